# Applies updated simulation results (case with 380 kV done) to Sheet1.
# Columns B, D, E, F, G, I, J, O for rows 2-25 receive new numeric values;
# columns A, C, H, K, L, M, N and the header row are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25)

$Bvals = @(1.219803940679071,1.069997499266549,0.9775436500530645,0.9397517899408854,0.9334695415156489,0.9770344427063833,1.168250105383493,1.539389405566169,1.809644099351658,1.932050532542917,1.978324344473322,1.968361989860625,1.935859098657602,1.915939811075305,1.801633628443199,1.731372246120657,1.690909632147168,1.677201143001923,1.738856893516015,1.945408141631617,2.079940855282643,2.008181005422728,1.735473297350438,1.439406031741498)
$Dvals = @(0.1740869811143284,0.1738283514950183,0.1737352727050379,0.1737138821688973,0.1737113296674551,0.1737349172359117,0.1739841729763967,0.1749940350993739,0.1760534290396052,0.1766042903173286,0.1768227941903859,0.1767752950894987,0.1766220683052069,0.1765295021117765,0.1760188154244986,0.1757231760749818,0.1755596213569888,0.1755053593357232,0.1757539758419355,0.1766668059812773,0.1773211216162025,0.1769666215723475,0.1757400312936568,0.1746650798694702)
$Evals = @(0.1584551540436152,0.1560976871575157,0.1547365595057464,0.1542035670625133,0.1541163722461825,0.1547292836604051,0.1576243489920373,0.1639891454149733,0.169088474725001,0.1715010782485393,0.1724280870050521,0.1722278419399359,0.1715770748056329,0.1711802087939134,0.1689326789279946,0.167577719709449,0.1668071234688995,0.166547712491564,0.167721052465474,0.1717678565883389,0.1744908541518555,0.1730303690042234,0.1676562256042899,0.1621932848350376)
$Fvals = @(1.586058537684849,1.537274724656484,1.508292599236057,1.496725533069323,1.494819505439395,1.508135617428664,1.569035736272298,1.696215314610981,1.794459758914059,1.840215008519493,1.857695457956879,1.853923867927108,1.841650044738657,1.834152050347342,1.791491030024417,1.765592947740373,1.750797082196556,1.745804605181775,1.768339483154278,1.845250978317722,1.896414790630246,1.869025219674057,1.76709748398801,1.660971659900625)
$Gvals = @(0.002415559989362569,0.002420193372685223,0.002423190390540851,0.002424450076845931,0.002424661568431892,0.002423207223543514,0.002417126089359947,0.002406401938415241,0.002399246738317698,0.002396147058881212,0.00239499548138893,0.002395242508514029,0.00239605187356056,0.002396550521006046,0.00239945242344497,0.002401272324867566,0.002402333704617743,0.00240269558458022,0.002401077080864969,0.002395813541770542,0.002392502880666877,0.002394258045721026,0.002401165303691459,0.002409175396314761)
$Ivals = @(3.961193417866525,3.481483530427994,3.186307868167489,3.065876933603903,3.04587117873993,3.184684259858585,3.795928923952232,4.989014804984947,5.861568061238756,6.25753114434093,6.407322560669854,6.375069214559744,6.269857673435979,6.205392518134659,5.835670441680065,5.608601612490133,5.477907426105048,5.433641513037287,5.632782903638144,6.30076505272973,6.736443306619492,6.503998686294949,5.621851005971507,4.666914883151151)
$Jvals = @(0.1913305446093858,0.1859881629284246,0.1828259976756428,0.1815669644110329,0.1813596857490793,0.1828088983185268,0.1894639112100833,0.2034575054670285,0.2143235631693869,0.2193960977187146,0.221335716487971,0.2209171483289651,0.2195552945535866,0.2187235669912297,0.2139946788174569,0.2111269263642015,0.2094896528478927,0.2089373882823509,0.2114309417489437,0.2199547937791095,0.225635036445496,0.2225933271110279,0.2112934607855692,0.1995698711105831)
$Ovals = @(3.975293305671926,3.829925653828582,3.743161388050623,3.708428181857073,3.702698365185597,3.742690441195975,3.92465161631327,4.301396925444635,4.590573542432821,4.724870106639514,4.776123705914472,4.765067556405256,4.729078768254283,4.707086574762002,4.581852524519547,4.505731717013816,4.462207568522672,4.447515346699447,4.513808127811615,4.739638709769963,4.889555756343043,4.809328534290898,4.510156039135268,4.197320322167002)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 2).Value2  = $Bvals[$i]   # column B
    $ws.Cells.Item($r, 4).Value2  = $Dvals[$i]   # column D
    $ws.Cells.Item($r, 5).Value2  = $Evals[$i]   # column E
    $ws.Cells.Item($r, 6).Value2  = $Fvals[$i]   # column F
    $ws.Cells.Item($r, 7).Value2  = $Gvals[$i]   # column G
    $ws.Cells.Item($r, 9).Value2  = $Ivals[$i]   # column I
    $ws.Cells.Item($r, 10).Value2 = $Jvals[$i]   # column J
    $ws.Cells.Item($r, 15).Value2 = $Ovals[$i]   # column O
}
